$d = $word.ActiveDocument

# 1. Title: "Updated Project Plan – Software App Wizards (SAW)"
#    -> "Updated Project Plan for Module 6 – Software App Wizards (SAW)"
$d.Content.Find.Execute(
    "Updated Project Plan – Software App Wizards (SAW)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Updated Project Plan for Module 6 – Software App Wizards (SAW)", 2
) | Out-Null

# 2. Programming Approach paragraph: merge the trailing "or just for the practice itself."
#    run(s) and append the new collaboration/extra-credit sentences.
$d.Content.Find.Execute(
    "or just for the practice itself.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "or just for the practice itself. As of this week we have also added another weekly collaboration session to work on the modules that we assigned each other last check in. Team leader has also reached out confirming the potential for extra credit by coding the payroll program in python. As such we will be hosting a coding session via discord to collaborate and code the payroll program. ",
    2
) | Out-Null

# 3. Sub-step owner lines: replace the space before the trailing name with a hyphen.
$d.Content.Find.Execute(
    "Initialize constants and income tax (Housecleaning) Alexia Erkman", $true, $false, $false, $false, $false,
    $true, 1, $false, "Initialize constants and income tax (Housecleaning)-Alexia Erkman", 2
) | Out-Null

$d.Content.Find.Execute(
    "Get employee input (hours and rate) Samuel Steinhardt", $true, $false, $false, $false, $false,
    $true, 1, $false, "Get employee input (hours and rate)-Samuel Steinhardt", 2
) | Out-Null

$d.Content.Find.Execute(
    "Run calculations (gross pay, deductions, net pay) Jayden Johnson", $true, $false, $false, $false, $false,
    $true, 1, $false, "Run calculations (gross pay, deductions, net pay)-Jayden Johnson", 2
) | Out-Null

$d.Content.Find.Execute(
    "Output payroll details for each employee Wai Moo", $true, $false, $false, $false, $false,
    $true, 1, $false, "Output payroll details for each employee-Wai Moo", 2
) | Out-Null

$d.Content.Find.Execute(
    "End of program with a final “End of Job” message Alexia Erkman", $true, $false, $false, $false, $false,
    $true, 1, $false, "End of program with a final “End of Job” message-Alexia Erkman", 2
) | Out-Null

# 4. Modules line: collapse the separately-run "(Input, Calculation/Run, Output), and
#    Termination(End of Job)" fragments into the preceding run's text.
$d.Content.Find.Execute(
    "(Input, Calculation/Run, Output), and Termination(End of Job)", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Input, Calculation/Run, Output), and Termination(End of Job)", 2
) | Out-Null
